# data: update experimental data
# Strip the trailing "-1" suffix from the benchmark/case names in column B
# of every results table, on every worksheet (SHANNON, SpaceEx, HyComp).

$wb = $excel.ActiveWorkbook

$names = @("altitude-display","altitude-display-int","ADC-bug-int","ADC-bug-d-int","car-controller","csma-aut","fisher-aut","hddi","water-tank","learning-factory","medical-monitor")

# SHANNON has five stacked 11-row tables starting at these rows.
$shannon = $wb.Worksheets.Item("SHANNON")
$shannonStarts = @(2, 15, 28, 41, 54)
foreach ($start in $shannonStarts) {
    for ($i = 0; $i -lt $names.Count; $i++) {
        $shannon.Cells.Item($start + $i, 2).Value = $names[$i]
    }
}

# SpaceEx and HyComp each have a single 11-row table starting at row 2.
$spaceEx = $wb.Worksheets.Item("SpaceEx")
for ($i = 0; $i -lt $names.Count; $i++) {
    $spaceEx.Cells.Item(2 + $i, 2).Value = $names[$i]
}

$hyComp = $wb.Worksheets.Item("HyComp")
for ($i = 0; $i -lt $names.Count; $i++) {
    $hyComp.Cells.Item(2 + $i, 2).Value = $names[$i]
}

# Reflect the updated active-sheet / selection state: SHANNON becomes the
# active tab, scrolled/selected near the bottom of the first table block.
$shannon.Activate()
$shannon.Application.ActiveWindow.ScrollRow = 29
$shannon.Range("B64").Select()

# SpaceEx's remembered selection moves to B12.
$spaceEx.Range("B12").Select()

# HyComp's remembered selection moves to B12 (and it's no longer the tab
# shown when the workbook re-opens, since SHANNON is now active).
$hyComp.Range("B12").Select()

$shannon.Activate()
